# ListObject#advanced_filter with help of limit
#
# Adds more sample rows below the existing "table3" ListObject on the
# "Tabelle1" sheet, grows the ListObject to cover them, turns off its
# AutoFilter dropdowns, and records the hidden _FilterDatabase defined
# name (localSheetId=2 -> Tabelle1) that Excel stores once a filter /
# advanced-filter has been applied against that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- new data rows appended below the current table body (D10:H16) ---

$ws.Range("D10").Value = 3
$ws.Range("E10").Value = "Eiffel"
$ws.Range("F10").Value = 50
$ws.Range("G10").Value = 0.5
$ws.Range("G10").NumberFormat = "h:mm"
$ws.Range("H10").Value = 30

$ws.Range("D11").Value = 3
$ws.Range("E11").Value = "Berta"
$ws.Range("G11").Value = 0.54166666666666663
$ws.Range("G11").NumberFormat = "h:mm"
$ws.Range("H11").Value = 40

# row 13 only carries the time number format, no content
$ws.Range("G13").NumberFormat = "h:mm"

$ws.Range("D14").Value = 3
$ws.Range("E14").Value = "Martha"

$ws.Range("D15").Value = 3
$ws.Range("E15").Value = "Paul"
$ws.Range("F15").Value = 40
$ws.Range("G15").Value = 0.5
$ws.Range("G15").NumberFormat = "h:mm"
$ws.Range("H15").Value = 80

$ws.Range("D16").Value = 1
$ws.Range("E16").Value = "Napoli"
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 0.41666666666666669
$ws.Range("G16").NumberFormat = "h:mm"
$ws.Range("H16").Value = 70

# --- grow the ListObject to the new extent and drop its AutoFilter UI ---

$lo = $ws.ListObjects.Item("table3")
$lo.Resize($ws.Range("D3:H16"))
$lo.ShowAutoFilter = $false

# --- hidden workbook-level _FilterDatabase name scoped to Tabelle1 ---

$name = $ws.Names.Add("_xlnm._FilterDatabase", "=Tabelle1!`$D`$3:`$H`$16")
$name.Visible = $false

# --- match the saved selection from the authoring session ---

[void]$ws.Activate()
[void]$ws.Range("J15").Select()
